$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1): F2 532 -> 533, F3 6370 -> 6373
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 533
$ws1.Range("F3").Value = 6373

# Update "全部类型" sheet (sheet4): F2 532 -> 533, F3 6370 -> 6373
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 533
$ws4.Range("F3").Value = 6373
